# Update computed results for Case 2_51 (380 kV) in pl_mw.xlsx, rows 2-25, columns B,C,E,F,G,H,I,K,L,N
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "E", "F", "G", "H", "I", "K", "L", "N")

$newValues = @{
    2 = @(0.6628252533777754, 0.05636454526420209, 0.08405757332559816, 0.4443680307746263, 1.016663603098138, 1.041725090818716, 1.005509539630399, 0.424796154726323, 0.2180205523869745, 1.911934984869209)
    3 = @(0.6235849028883536, 0.05399335241401104, 0.08339352120374954, 0.387822817061874, 1.020924513557887, 1.048347776992102, 1.012580043170352, 0.388305896346651, 0.2108835005249006, 1.932299460346287)
    4 = @(0.5997605487285966, 0.0525176477162006, 0.08302523262647554, 0.3531389305168915, 1.024149801122419, 1.05285443413851, 1.017404798925682, 0.3660349307352817, 0.2066213008989735, 1.94543939277232)
    5 = @(0.5901201821920381, 0.05191130000524424, 0.08288509085488904, 0.3390132514313251, 1.025617162556422, 1.054801691921583, 1.019492478761617, 0.3569934281344729, 0.2049146490264775, 1.95095394709097)
    6 = @(0.5885235440301813, 0.05181031494933563, 0.08286242131253907, 0.336668177824194, 1.02587005582469, 1.055131723356823, 1.019846477567835, 0.3554941632393991, 0.2046330883012928, 1.951879293246941)
    7 = @(0.5996302583457691, 0.05250949050240195, 0.08302330235582289, 0.3529483938344953, 1.024168971024494, 1.052880247028156, 1.017432461960478, 0.3659128553474602, 0.2065981619189472, 1.945513116474347)
    8 = @(0.6492395652282426, 0.05555105632500101, 0.08382043399625871, 0.4248636149813336, 1.01800627501261, 1.0439172326099, 1.007847104051979, 0.4121866588764931, 0.2155348281117568, 1.918824632469081)
    9 = @(0.7486435632488337, 0.06135935063794307, 0.08569582341913318, 0.5661985755041457, 1.010760509782884, 1.029832870115868, 0.9928869800629627, 0.5039836603803565, 0.2340101386414375, 1.87153545076783)
    10 = @(0.822953323261288, 0.06553287215865566, 0.08726327173333814, 0.6702781546542269, 1.008397564118781, 1.02161236890035, 0.9842364537948711, 0.5720624383335178, 0.2481631999587819, 1.839867594846749)
    11 = @(0.8570335157829732, 0.06741145953991179, 0.08801737477945437, 0.7176906081379002, 1.007967780790409, 1.018334290401683, 0.9808097856899636, 0.6031701545722683, 0.2547276317861531, 1.826128638299531)
    12 = @(0.8699781253475578, 0.06811997771788469, 0.08830882177360877, 0.7356546913071611, 1.007897972551575, 1.017159303219444, 0.9795853510914654, 0.6149694759843669, 0.2572315149997593, 1.821021982728869)
    13 = @(0.8671885383119502, 0.06796751279441082, 0.08824579205536409, 0.7317853510981394, 1.007908870668317, 1.01740940692595, 0.9798457999431065, 0.6124274197237582, 0.2566914560409828, 1.82211752179408)
    14 = @(0.8580976945362693, 0.06746980706385841, 0.08804123448033252, 0.7191683204515869, 1.007960173954004, 1.018236293792867, 0.9807075839974146, 0.6041405031744489, 0.2549332658668817, 1.825706586390728)
    15 = @(0.8525343762622697, 0.06716457548058941, 0.08791670279786246, 0.7114413442032514, 1.008003707521382, 1.018751426361391, 0.9812449819927096, 0.599067056026314, 0.2538586757033698, 1.82791749537721)
    16 = @(0.8207316153344664, 0.06540970103819177, 0.08721481388542784, 0.6671810134426437, 1.008438647892234, 1.021835884210788, 0.9844706310690228, 0.5700322264317208, 0.2477367316913899, 1.840778895371604)
    17 = @(0.8012920277897706, 0.06432803465352777, 0.08679473056313114, 0.6400460337125793, 1.008870826017869, 1.023846285984263, 0.9865797356599018, 0.5522554542966986, 0.2440133787261516, 1.848839821094575)
    18 = @(0.7901369318084903, 0.06370400867225356, 0.08655697599903434, 0.6244449056556647, 1.00918011605863, 1.025046051984788, 0.9878406962825039, 0.5420437694558586, 0.2418836805917692, 1.853539019590877)
    19 = @(0.7863644973311636, 0.06349240100600184, 0.08647714103836179, 0.619163680173358, 1.009295258564904, 1.025459732606706, 0.9882758543205981, 0.5385885206428895, 0.241164643011345, 1.855140865158756)
    20 = @(0.8033587139699137, 0.06444337440139947, 0.08683904915434582, 0.642933953830422, 1.008818535143206, 1.023627779875582, 0.9863502642506035, 0.5541464748158944, 0.2444085073353079, 1.847975224985132)
    21 = @(0.8607668367791632, 0.067616072903391, 0.08810115845555089, 0.7228739723491628, 1.007942581160236, 1.017991616420261, 0.9804524709208593, 0.6065740451950319, 0.255449198879532, 1.82464978528224)
    22 = @(0.8985144802032323, 0.06967293664935426, 0.08896031001271609, 0.7751780083420101, 1.007911908134957, 1.014694792992159, 0.9770244372158388, 0.6409519994677453, 0.262770256818925, 1.809964803227519)
    23 = @(0.878347162347211, 0.06857667271916057, 0.08849863380405765, 0.7472568307830727, 1.007878646001302, 1.016418986033344, 0.9788150009347945, 0.6225935808053578, 0.2588532536242525, 1.817751226467383)
    24 = @(0.8024242997498163, 0.06439123603625063, 0.08681900101563045, 0.6416283278902171, 1.008841986392397, 1.02372642955288, 0.9864538574873265, 0.5532915181887859, 0.2442298357133268, 1.848365906836634)
    25 = @(0.7215267032461554, 0.05980463188654284, 0.0851551320170536, 0.5279251897347166, 1.012201534414203, 1.033269379746841, 0.9965231682644387, 0.4790380192473833, 0.2289103838387081, 1.883788263061899)
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "{0}{1}" -f $columns[$i], $row
        $ws.Range($cellRef).Value = $rowValues[$i]
    }
}

